$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.052.67'
$ws.Range('E2').Value = '  +0.69%  '

# Row 3
$ws.Range('D3').Value = '1.683.81'
$ws.Range('E3').Value = '  +1.00%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.27%  '

# Row 6
$ws.Range('E6').Value = '  -2.14%  '

# Row 7
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.254'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.29%  '

# Row 9
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '21.50'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.43%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0623'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.67%  '

# Row 11
$ws.Range('E11').Value = '  -0.61%  '

# Row 12
$ws.Range('E12').Value = '  +0.94%  '

# Row 13
$ws.Range('D13').Value = '1.680.57'
$ws.Range('E13').Value = '  +0.18%  '

# Row 14
$ws.Range('E14').Value = '  +0.67%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.534'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.92%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.41'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.78%  '

# Row 17
$ws.Range('D17').Value = '27.059.24'
$ws.Range('E17').Value = '  +0.68%  '

# Row 18
$ws.Range('E18').Value = '  +5.32%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '236.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.23%  '

# Row 20
$ws.Range('E20').Value = '  +0.71%  '

# Row 22
$ws.Range('E22').Value = '  +0.40%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.16%  '

# Row 24
$ws.Range('E24').Value = '  -3.68%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.07'
$ws.Range('D25').Style = 'Normal'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.48%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.73%  '

# Row 28
$ws.Range('E28').Value = '  -2.37%  '

# Row 29
$ws.Range('E29').Value = '  +0.07%  '

# Row 30
$ws.Range('E30').Value = '  +0.55%  '

# Row 31
$ws.Range('E31').Value = '  -0.02%  '

# Row 32
$ws.Range('E32').Value = '  +0.52%  '

# Row 33
$ws.Range('D33').Value = '1.526.61'
$ws.Range('E33').Value = '  +4.30%  '

# Row 34
$ws.Range('E34').Value = '  +0.98%  '

# Row 35
$ws.Range('E35').Value = '  +4.74%  '

# Row 36
$ws.Range('E36').Value = '  -0.48%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.590'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.47%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.920'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.66%  '

# Row 39
$ws.Range('E39').Value = '  +3.67%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.74%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.94%  '

# Row 42
$ws.Range('E42').Value = '  -0.08%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.80%  '

# Row 44
$ws.Range('E44').Value = '  -0.70%  '

# Row 45
$ws.Range('D45').Value = '1.825.42'
$ws.Range('E45').Value = '  +0.56%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.781'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.08%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.18%  '

# Row 48
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.104'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.24%  '

# Row 49
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.02%  '

# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.36%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0507'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.18%  '
